# Update countries & provincias Spain
# Applies the COVID data refresh captured in the diff:
#  - Argentina's row moves above Sudafrica (its totals overtook Sudafrica's)
#  - Several country rows get refreshed case/death counters
#  - The "datos actualizados" timestamp string is bumped

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 586377
$ws.Range("C4").Value = 26077
$ws.Range("D4").Value = 36218
$ws.Range("E4").Value = 526549
$ws.Range("G4").Value = 1505
$ws.Range("H4").Value = 23610

# Row 44 - Panama
$ws.Range("B44").Value = 3472
$ws.Range("C44").Value = 72
$ws.Range("D44").Value = 61
$ws.Range("E44").Value = 3317
$ws.Range("G44").Value = 7
$ws.Range("H44").Value = 94

# Row 54 now becomes Argentina (its totals passed Sudafrica's), row 55 becomes Sudafrica
$ws.Range("A54").Value = "Argentina"
$ws.Range("B54").Value = 2277
$ws.Range("C54").Value = 135
$ws.Range("D54").Value = 515
$ws.Range("E54").Value = 1664
$ws.Range("F54").Value = 83
$ws.Range("G54").Value = 8
$ws.Range("H54").Value = 98

$ws.Range("A55").Value = "Sudafrica"
$ws.Range("B55").Value = 2272
$ws.Range("C55").Value = 99
$ws.Range("D55").Value = 410
$ws.Range("E55").Value = 1835
$ws.Range("F55").Value = 7
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 27

# Row 95 - Burkina Faso
$ws.Range("B95").Value = 515
$ws.Range("C95").Value = 18
$ws.Range("E95").Value = 327

# Row 104 - San Marino
$ws.Range("B104").Value = 371
$ws.Range("C104").Value = 15
$ws.Range("E104").Value = 282
$ws.Range("F104").Value = 15
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 36

# Row 125 - Guadalupe
$ws.Range("B125").Value = 145
$ws.Range("C125").Value = 2
$ws.Range("E125").Value = 70

# Row 128 - Gibraltar
$ws.Range("D128").Value = 93
$ws.Range("E128").Value = 36

# Updated "last refreshed" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 01:52"
